$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column B (original file name) and C (durations) for rows 5-8
$ws.Range("B5").Value = "rfinnegan_tape 1_side 1"
$ws.Range("C5").Value = 0.16527777777777777

$ws.Range("B6").Value = "rfinnegan_tape 1_side 1"
$ws.Range("C6").Value = 0.10277777777777779

$ws.Range("B7").Value = "rfinnegan_tape 1_side 1"
$ws.Range("C7").Value = 0.33263888888888887

$ws.Range("B8").Value = "rfinnegan_tape 1_side 1"
$ws.Range("C8").Value = 0.14930555555555555

# Apply the same duration number format style as existing C column cells
$ws.Range("C5:C8").NumberFormat = $ws.Range("C2").NumberFormat

# Fill in column D (new file name) top-to-bottom first
$ws.Range("D5").Value = "tape1_side1_tale4"
$ws.Range("D6").Value = "tape1_side1_tale5"
$ws.Range("D7").Value = "tape1_side1_tale6"
$ws.Range("D8").Value = "tape1_side1_tale7"

# Then fill in column E (comments) top-to-bottom
$ws.Range("E6").Value = "SONG"
$ws.Range("E7").Value = "maybe several tales by the same speaker"
$ws.Range("E8").Value = "song in the end"

# Update the active selection to reflect the final cursor position
$ws.Range("K17").Select()
